$wb = $excel.ActiveWorkbook

# This workbook tracks localization handback status for three files
# (189daec5*.png, 3c6023d3*.png, 6d7b2b2c*.md, be1638a0*.md, bf579bfb*.md)
# across an "Overview" sheet and one detail sheet per target locale
# (zh-cn, de-de). Re-generating the handback report refreshes the
# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the file that was just
# processed (be1638a0-53d9-4644-b3a3-2a2b1a97907b.md), row 5 on every
# sheet, leaving every other row untouched.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G5").Value = "2016-11-14 05:48:59"

# zh-cn detail sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
$wsZhCn.Range("H5").Value = "2016-11-14 05:48:46"
$wsZhCn.Range("K5").Value = "2016-11-14 05:49:39"

# de-de detail sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
$wsDeDe.Range("H5").Value = "2016-11-14 05:48:59"
$wsDeDe.Range("K5").Value = "2016-11-14 05:49:57"
